$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row added for "Change Vcenter Details" test entry (waveEdit.py)
$ws.Range("B32").Value = "./TestData/OneForAll/editVCenterData.xlsx"
$ws.Range("A32").Value = "Change Vcenter Details"
$ws.Range("C32").Value = "NA"
$ws.Range("D32").Value = "NA"

# Match the author's final selection/scroll position on save
$ws.Range("A32").Select() | Out-Null
